# Update "想去人数" (want-to-go count) values in column F
# for the "展览" (Exhibition) sheet and the "全部类型" (All types) sheet.
# Matching rows are identified by row number verified against both sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Row -> new F value for "展览" sheet
$sheet1Updates = @{
    2  = 15189
    3  = 19598
    5  = 179
    14 = 221
    15 = 255
    16 = 77
    17 = 1530
    20 = 119
    21 = 250
    22 = 8262
    26 = 71
    27 = 1277
    28 = 27
    31 = 6565
    32 = 137
    34 = 192
    35 = 159
    36 = 313
    37 = 5621
    38 = 1020
    39 = 33
}

foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# Row -> new F value for "全部类型" sheet
$sheet4Updates = @{
    2  = 15189
    3  = 19598
    5  = 179
    14 = 221
    15 = 255
    16 = 77
    17 = 1530
    21 = 119
    22 = 250
    23 = 8262
    27 = 71
    28 = 1277
    29 = 27
    34 = 6565
    35 = 137
    37 = 192
    38 = 159
    39 = 313
    40 = 5621
    41 = 1020
    42 = 33
}

foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
